# "Refined metadata to be additional tab"
#
# 1. Update the F-column (time_taken) timestamps on the "data" sheet.
# 2. Add a new "metadata" worksheet after "data" with a header row and one
#    data row describing the PanelApp query that produced the data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("F2").Value = "2021-10-05 14:19:41.281548"
$ws.Range("F3").Value = "2021-10-05 14:19:41.281556"
$ws.Range("F4").Value = "2021-10-05 14:19:41.281559"
$ws.Range("F5").Value = "2021-10-05 14:19:41.281562"
$ws.Range("F6").Value = "2021-10-05 14:19:41.281565"
$ws.Range("F7").Value = "2021-10-05 14:19:41.281568"
$ws.Range("F8").Value = "2021-10-05 14:19:41.281570"
$ws.Range("F9").Value = "2021-10-05 14:19:41.281573"
$ws.Range("F10").Value = "2021-10-05 14:19:41.281575"
$ws.Range("F11").Value = "2021-10-05 14:19:41.281578"
$ws.Range("F12").Value = "2021-10-05 14:19:41.281581"
$ws.Range("F13").Value = "2021-10-05 14:19:41.281583"
$ws.Range("F14").Value = "2021-10-05 14:19:41.281586"
$ws.Range("F15").Value = "2021-10-05 14:19:41.281588"
$ws.Range("F16").Value = "2021-10-05 14:19:41.281591"
$ws.Range("F17").Value = "2021-10-05 14:19:41.281593"
$ws.Range("F18").Value = "2021-10-05 14:19:41.281596"
$ws.Range("F19").Value = "2021-10-05 14:19:41.281599"
$ws.Range("F20").Value = "2021-10-05 14:19:41.281601"
$ws.Range("F21").Value = "2021-10-05 14:19:41.281604"
$ws.Range("F22").Value = "2021-10-05 14:19:41.281606"

# Add the "metadata" worksheet right after "data".
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row (B1:G1) and the A2 row-index cell reuse the exact same style
# as the "data" sheet's header row (bold, centered, thin-bordered) by
# copying the format across, same as how the rest of the workbook already
# shares that single style.
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"
$meta.Range("A2").Value = 0

# Data row describing the panel query.
$meta.Range("B2").Value = "Congenital adrenal hypoplasia"
$meta.Range("C2").Value = 145
$meta.Range("E2").Value = "2021-01-28T13:04:07.145691Z"
$meta.Range("F2").Value = "2021-10-05 14:19:41.277966"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/145/?format=json"

# D2 ("2.6") must stay a literal text value, not be coerced to the number
# 2.6 -- route it through a text formula + paste-special-values so it lands
# as a plain (unstyled) string cell, same as every other sheet-provided
# string on these sheets.
$meta.Range("Z1").Formula = "=TEXT(2.6,""0.0"")"
$meta.Range("Z1").Copy()
$meta.Range("D2").PasteSpecial(-4163)
$meta.Range("Z1").Clear()

$meta.Range("A1").Select() | Out-Null

# Keep "data" as the active/selected sheet, as it was before this edit.
$ws.Activate() | Out-Null
